{"js": "// The document has a series of bold/italic \"heading\" paragraphs (e.g.\n// \"CREATION OF FUTURE:\") introducing each section. This adds a new one,\n// \"EXECUTERS:\", right after the paragraph that ends with \"...our code\n// will start its execution and then termination.\" and right before the\n// paragraph that begins the \"In async rust we have executers...\" section.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst marker = \"our code will start its execution and then termination.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the paragraph ending in: \" + marker);\n}\n\n// Insert the new heading paragraph immediately after the target paragraph.\nconst newPara = target.insertParagraph(\"EXECUTERS:\", Word.InsertLocation.after);\n\n// Match the formatting used by the document's other section headings:\n// Times New Roman, bold, italic, 12pt (OOXML sz/szCs = 24 half-points).\nconst headingFont = { name: \"Times New Roman\", bold: true, italic: true, size: 12 };\nnewPara.font.set(headingFont);\n// Also stamp the paragraph-mark run properties (the formatting that would\n// apply to text typed at the very end of the paragraph) so it matches the\n// sibling headings exactly.\nnewPara.getRange(\"End\").font.set(headingFont);\n\nawait context.sync();\n", "ps1": "# The document has a series of bold/italic \"heading\" paragraphs (e.g.\n# \"CREATION OF FUTURE:\") introducing each section. This adds a new one,\n# \"EXECUTERS:\", right after the paragraph that ends with \"...our code\n# will start its execution and then termination.\" and right before the\n# paragraph that begins the \"In async rust we have executers...\" section.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"our code will start its execution and then termination.\")\n\nif (-not $found) {\n    throw \"Could not locate the target sentence in the document.\"\n}\n\n# $rng now spans the found sentence; its paragraph is the insertion anchor.\n$targetPara = $rng.Paragraphs(1)\n\n$anchorRange = $targetPara.Range\n$anchorRange.Collapse(0)   # wdCollapseEnd\n$anchorRange.InsertParagraphAfter()\n\n# The freshly inserted paragraph immediately follows the target paragraph.\n$newPara = $targetPara.Next()\n$newRange = $newPara.Range\n$newRange.Text = \"EXECUTERS:\"\n\n# Match the formatting used by the document's other section headings:\n# Times New Roman, bold, italic, 12pt.\n$newRange.Font.Name = \"Times New Roman\"\n$newRange.Font.Bold = 1\n$newRange.Font.Italic = 1\n$newRange.Font.Size = 12\n"}
